$wb = $excel.ActiveWorkbook

# --- Planilha1 sheet (dimension A1:L6) ---
$ws1 = $wb.Worksheets.Item("Planilha1")

# Row 6, column A currently shares the same string as A2 ("ISAQUEab").
# Update A6's text to a new value; this changes the shared string that A6
# references (which A2 used to reference too, but A2 gets its own new value below).
$ws1.Range("A6").Value = "ISAQUEhg715"

# Update A2 to a brand-new value, so it receives its own shared-string entry.
$ws1.Range("A2").Value = "ISAQUEz1"

# Update the selection on Planilha1 to A2.
$ws1.Select()
$ws1.Range("A2").Select()

# --- Planilha2 sheet (dimension A1:A5) ---
$ws2 = $wb.Worksheets.Item("Planilha2")
$ws2.Select()
$ws2.Range("B5").Select()

# Reselect Planilha1 so it stays the active/visible tab as before.
$ws1.Select()
$ws1.Range("A2").Select()
